$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the old "_GoBack" bookmark that sits right after the sentence
#    "...at the end of the internship."
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 2) In the "Student:" paragraph (the one introducing the professional
#    objectives list), add a space after the colon: ":" -> ": "
# ---------------------------------------------------------------------------
$br = [char]11
$full = $d.Content.Text
$anchor = "Student:" + $br + "My professional objectives include learning"
$anchorIdx = $full.IndexOf($anchor)
if ($anchorIdx -lt 0) {
    throw "Could not locate the 'Student:' professional-objectives paragraph"
}
$colonIdx = $anchorIdx + ("Student".Length)
$colonRange = $d.Range($colonIdx, $colonIdx + 1)
if ($colonRange.Text -ne ":") {
    throw "Unexpected text at colon position: [$($colonRange.Text)]"
}
$colonRange.Text = ": "

# ---------------------------------------------------------------------------
# 3) Re-insert the "_GoBack" bookmark in the middle of "learning", right
#    after "learni" (the bookmark is an empty/collapsed range).
# ---------------------------------------------------------------------------
$full = $d.Content.Text
$learningIdx = $full.IndexOf("learning more about the new office")
if ($learningIdx -lt 0) {
    throw "Could not locate 'learning more about the new office'"
}
$bookmarkPos = $learningIdx + ("learni").Length
$bookmarkRange = $d.Range($bookmarkPos, $bookmarkPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

# ---------------------------------------------------------------------------
# 4) "new technology, gain upward feedback from team members " ->
#    "new technology, gaining upward feedback from team members "
# ---------------------------------------------------------------------------
$found = $d.Content.Find.Execute(
    "new technology, gain upward feedback from team members ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "new technology, gaining upward feedback from team members ", 2)
if (-not $found) {
    throw "Could not find/replace the 'new technology, gain' sentence"
}

# ---------------------------------------------------------------------------
# 5) The long "Strong"-styled sentence at the end of the paragraph:
#      "to establish further what can be improved upon, Improve presentation
#       skills,  gain a more in-depth perspective of how departments within
#       the company are run, network, and increase my performance metrics."
#    becomes:
#      "to establish further what can be improved upon, improving my
#       presentation skills, gaining a more in-depth perspective of how
#       departments within the company are run, network, and increasing my
#       performance metrics."
#    and loses its explicit dark-grey font color (w:color val="30393F").
# ---------------------------------------------------------------------------
$nbsp = [char]160
$oldSentence = "to establish further what can be improved upon, Improve presentation skills, " + $nbsp + "gain a more in-depth perspective of how departments within the company are run, network, and increase my performance metrics."
$newSentence = "to establish further what can be improved upon, improving my presentation skills, gaining a more in-depth perspective of how departments within the company are run, network, and increasing my performance metrics."

$full = $d.Content.Text
$sentIdx = $full.IndexOf($oldSentence)
if ($sentIdx -lt 0) {
    throw "Could not locate the long 'to establish further...' sentence"
}
$sentRange = $d.Range($sentIdx, $sentIdx + $oldSentence.Length)
if ($sentRange.Text -ne $oldSentence) {
    throw "Sentence range text mismatch: [$($sentRange.Text)]"
}
$sentRange.Text = $newSentence
# Re-acquire the range after the text assignment (length changed) and strip
# the explicit font color so it reverts to the style's default.
$sentRange2 = $d.Range($sentIdx, $sentIdx + $newSentence.Length)
$sentRange2.Font.Color = -16777216  ## wdColorAutomatic
